$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.017.08'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.314.96'
$ws.Range('E3').Value = '  +2.39%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''253.49'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').Value = '''0.632'
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('D7').Value = '''75.98'
$ws.Range('E7').Value = '  +7.41%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.642'
$ws.Range('E9').Value = '  -4.26%  '
$ws.Range('D10').Value = '''39.55'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('E11').Value = '  +1.55%  '
$ws.Range('E12').Value = '  -1.72%  '
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('D14').Value = '2.661.67'
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('D15').Value = '''15.45'
$ws.Range('E15').Value = '  +3.93%  '
$ws.Range('D16').Value = '''0.883'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').Value = '2.321.44'
$ws.Range('E17').Value = '  +2.51%  '
$ws.Range('D18').Value = '42.968.80'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('E19').Value = '  +2.88%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').Value = '''72.94'
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('D22').Value = '''236.89'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('E23').Value = '  +5.53%  '
$ws.Range('E24').Value = '  -0.55%  '
$ws.Range('D25').Value = '''11.66'
$ws.Range('E25').Value = '  -1.37%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '''2.42'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('D29').Value = '''21.32'
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').Value = '''167.32'
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('D31').Value = '''6.36'
$ws.Range('E31').Value = '  +1.34%  '
$ws.Range('D32').Value = '''0.0843'
$ws.Range('E32').Value = '  +8.98%  '
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').Value = '''30.43'
$ws.Range('E34').Value = '  +4.22%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = '''0.128'
$ws.Range('E35').Value = '  +1.69%  '
$ws.Range('D36').Value = '''4.59'
$ws.Range('E36').Value = '  +11.12%  '
$ws.Range('E37').Value = '  +3.43%  '
$ws.Range('E38').Value = '  -2.17%  '
$ws.Range('D39').Value = '''13.92'
$ws.Range('E39').Value = '  +14.52%  '
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('D41').Value = '''5.89'
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('E42').Value = '  +7.55%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''9.22'
$ws.Range('E43').Value = '  +3.18%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').Value = '''62.61'
$ws.Range('E44').Value = '  -3.07%  '
$ws.Range('D45').Value = '''4.91'
$ws.Range('E45').Value = '  -2.93%  '
$ws.Range('D46').Value = '''107.77'
$ws.Range('E46').Value = '  +13.82%  '
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('D51').Value = '''4.32'
$ws.Range('E51').Value = '  -1.55%  '
